$d = $word.ActiveDocument

# Locate the run of text that needs to be split/corrected.
$rng = $d.Content
$found = $rng.Find.Execute("W: D x H (H neurons, D dimensional layer x)")

if ($found) {
    # Expand to the whole paragraph (including its end-of-paragraph mark)
    # so InsertXML replaces the paragraph's runs in place while keeping
    # the paragraph's own properties (style, numbering, paraId, rsids...).
    # Using a freshly-obtained Range (rather than re-using the Find range
    # object) is what makes InsertXML merge into the existing paragraph.
    $para = $rng.Paragraphs(1)
    $pStart = $para.Range.Start
    $pEnd = $para.Range.End
    $rng = $d.Range($pStart, $pEnd)

    # Replace the single run with several runs carrying the corrected
    # text, matching the author's "Dec 2021 4)a) correction":
    #   "W: " / "H x D" / " (H" / " hidden" / " neurons, D dimensional layer x)" / " order matters!"
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' +
        '<w:r><w:t xml:space="preserve">W: </w:t></w:r>' +
        '<w:r><w:t>H x D</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> (H</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> hidden</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> neurons, D dimensional layer x)</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> order matters!</w:t></w:r>' +
        '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $rng.InsertXML($xml) | Out-Null
}
